$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C:G (5 contiguous columns), rows 2-25
$dataCG = New-Object 'object[,]' 24,5
$dataCG[0,0] = 7.553218057833775
$dataCG[0,1] = 8.487816639668129
$dataCG[0,2] = 12.08984770114114
$dataCG[0,3] = 42.74871832541115
$dataCG[0,4] = 3.731182238444394
$dataCG[1,0] = 7.516751532984305
$dataCG[1,1] = 8.528710039019513
$dataCG[1,2] = 12.05712619318125
$dataCG[1,3] = 42.61672848073314
$dataCG[1,4] = 3.735603616638222
$dataCG[2,0] = 7.495676701590162
$dataCG[2,1] = 8.555047374885385
$dataCG[2,2] = 12.03975304369503
$dataCG[2,3] = 42.54924594592943
$dataCG[2,4] = 3.738455142323995
$dataCG[3,0] = 7.487426128552534
$dataCG[3,1] = 8.56608922740991
$dataCG[3,2] = 12.03336049992926
$dataCG[3,3] = 42.52515703174032
$dataCG[3,4] = 3.739651702055775
$dataCG[4,0] = 7.486076705806541
$dataCG[4,1] = 8.567941398998602
$dataCG[4,2] = 12.03234061998205
$dataCG[4,3] = 42.5213630264074
$dataCG[4,4] = 3.739852480261507
$dataCG[5,0] = 7.495564055973208
$dataCG[5,1] = 8.555195037005918
$dataCG[5,2] = 12.03966404495347
$dataCG[5,3] = 42.54890726633077
$dataCG[5,4] = 3.738471139502274
$dataCG[6,0] = 7.540375470138487
$dataCG[6,1] = 8.501661817917693
$dataCG[6,2] = 12.0780036242926
$dataCG[6,3] = 42.70039592723249
$dataCG[6,4] = 3.732678428175805
$dataCG[7,0] = 7.638377660464083
$dataCG[7,1] = 8.406421078400518
$dataCG[7,2] = 12.17456209129029
$dataCG[7,3] = 43.10479835411721
$dataCG[7,4] = 3.722397612497694
$dataCG[8,0] = 7.716111852337705
$dataCG[8,1] = 8.342368416084263
$dataCG[8,2] = 12.25823672497785
$dataCG[8,3] = 43.46670058373367
$dataCG[8,4] = 3.715492592150843
$dataCG[9,0] = 7.752619234597866
$dataCG[9,1] = 8.314511973946571
$dataCG[9,2] = 12.29899055857364
$dataCG[9,3] = 43.64517344343926
$dataCG[9,4] = 3.712490096693397
$dataCG[10,0] = 7.766599486815084
$dataCG[10,1] = 8.304147635654877
$dataCG[10,2] = 12.31480225537773
$dataCG[10,3] = 43.71471996785385
$dataCG[10,4] = 3.711372910636965
$dataCG[11,0] = 7.763581810987152
$dataCG[11,1] = 8.306371584166998
$dataCG[11,2] = 12.31138019138794
$dataCG[11,3] = 43.69965507265297
$dataCG[11,4] = 3.711612638558992
$dataCG[12,0] = 7.753766330524788
$dataCG[12,1] = 8.313655599864106
$dataCG[12,2] = 12.30028384466995
$dataCG[12,3] = 43.65085591564321
$dataCG[12,4] = 3.71239778923642
$dataCG[13,0] = 7.747774069620752
$dataCG[13,1] = 8.318141271296147
$dataCG[13,2] = 12.29353614437225
$dataCG[13,3] = 43.62121973317802
$dataCG[13,4] = 3.712881290239451
$dataCG[14,0] = 7.713748294780211
$dataCG[14,1] = 8.344214687012425
$dataCG[14,2] = 12.25562684839021
$dataCG[14,3] = 43.45531322359215
$dataCG[14,4] = 3.715691589846709
$dataCG[15,0] = 7.693161657359639
$dataCG[15,1] = 8.360538094717576
$dataCG[15,2] = 12.23305436877426
$dataCG[15,3] = 43.35706181137156
$dataCG[15,4] = 3.717451024172523
$dataCG[16,0] = 7.681429164437336
$dataCG[16,1] = 8.370047498404459
$dataCG[16,2] = 12.22032495310592
$dataCG[16,3] = 43.30185498383792
$dataCG[16,4] = 3.71847606150176
$dataCG[17,0] = 7.677475633322143
$dataCG[17,1] = 8.37328793403803
$dataCG[17,2] = 12.21605878564608
$dataCG[17,3] = 43.2833877721908
$dataCG[17,4] = 3.718825368430512
$dataCG[18,0] = 7.695341986327024
$dataCG[18,1] = 8.358787957237784
$dataCG[18,2] = 12.23543104693445
$dataCG[18,3] = 43.36738598861324
$dataCG[18,4] = 3.717262379036129
$dataCG[19,0] = 7.756645222635463
$dataCG[19,1] = 8.311511105901111
$dataCG[19,2] = 12.3035328860499
$dataCG[19,3] = 43.66513638488195
$dataCG[19,4] = 3.712166635378074
$dataCG[20,0] = 7.79761362413586
$dataCG[20,1] = 8.281687179844383
$dataCG[20,2] = 12.3502471173422
$dataCG[20,3] = 43.87115795895689
$dataCG[20,4] = 3.7089515802298
$dataCG[21,0] = 7.775668418244086
$dataCG[21,1] = 8.297506461655386
$dataCG[21,2] = 12.32511571570237
$dataCG[21,3] = 43.76016526108278
$dataCG[21,4] = 3.710657011690709
$dataCG[22,0] = 7.694355937973583
$dataCG[22,1] = 8.359578805485267
$dataCG[22,2] = 12.2343557783096
$dataCG[22,3] = 43.36271444155352
$dataCG[22,4] = 3.717347623394977
$dataCG[23,0] = 7.610825582555254
$dataCG[23,1] = 8.431145459485622
$dataCG[23,2] = 12.1461802348498
$dataCG[23,3] = 42.98395412061892
$dataCG[23,4] = 3.725064327455917
$ws.Range("C2:G25").Value = $dataCG

# Column I, rows 2-25
$dataI = New-Object 'object[,]' 24,1
$dataI[0,0] = 36.85582234438661
$dataI[1,0] = 36.80845167361132
$dataI[2,0] = 36.79009882933705
$dataI[3,0] = 36.78531066597314
$dataI[4,0] = 36.78467782388999
$dataI[5,0] = 36.79002337219057
$dataI[6,0] = 36.83725718788548
$dataI[7,0] = 37.01529371850945
$dataI[8,0] = 37.19834501413722
$dataI[9,0] = 37.29294769196509
$dataI[10,0] = 37.33039465568663
$dataI[11,0] = 37.32225773147631
$dataI[12,0] = 37.29599599602799
$dataI[13,0] = 37.28012104610329
$dataI[14,0] = 37.19239003528993
$dataI[15,0] = 37.14146845599966
$dataI[16,0] = 37.11324686858443
$dataI[17,0] = 37.10387496332872
$dataI[18,0] = 37.14677873683762
$dataI[19,0] = 37.30366572216623
$dataI[20,0] = 37.41565552299684
$dataI[21,0] = 37.35502215685991
$dataI[22,0] = 37.14437467572475
$dataI[23,0] = 36.95795716698232
$ws.Range("I2:I25").Value = $dataI

# Columns K:L (2 contiguous columns), rows 2-25
$dataKL = New-Object 'object[,]' 24,2
$dataKL[0,0] = 21.62106131429388
$dataKL[0,1] = 9.908629492983968
$dataKL[1,0] = 21.23929134009351
$dataKL[1,1] = 9.91616875452278
$dataKL[2,0] = 21.00733586716919
$dataKL[2,1] = 9.922681130767241
$dataKL[3,0] = 20.91355990627451
$dataKL[3,1] = 9.925807709469817
$dataKL[4,0] = 20.89803747146302
$dataKL[4,1] = 9.926355400551373
$dataKL[5,0] = 21.00606797089776
$dataKL[5,1] = 9.92272138400163
$dataKL[6,0] = 21.48899642512682
$dataKL[6,1] = 9.91083770754878
$dataKL[7,0] = 22.44976457808055
$dataKL[7,1] = 9.902508791858528
$dataKL[8,0] = 23.15651310963117
$dataKL[8,1] = 9.905556133889835
$dataKL[9,0] = 23.47669755964242
$dataKL[9,1] = 9.908937885406074
$dataKL[10,0] = 23.5976327208596
$dataKL[10,1] = 9.910505510422372
$dataKL[11,0] = 23.57160291822204
$dataKL[11,1] = 9.910155129138081
$dataKL[12,0] = 23.48665389501808
$dataKL[12,1] = 9.909061102505975
$dataKL[13,0] = 23.43457604074062
$dataKL[13,1] = 9.90842835877932
$dataKL[14,0] = 23.13555119943113
$dataKL[14,1] = 9.905375293532707
$dataKL[15,0] = 22.95168272124364
$dataKL[15,1] = 9.90401354067489
$dataKL[16,0] = 22.84580925463607
$dataKL[16,1] = 9.903418128194314
$dataKL[17,0] = 22.80994600190505
$dataKL[17,1] = 9.90324878692717
$dataKL[18,0] = 22.97126886441942
$dataKL[18,1] = 9.9041390606297
$dataKL[19,0] = 23.51161488771192
$dataKL[19,1] = 9.909374654953155
$dataKL[20,0] = 23.86289265096975
$dataKL[20,1] = 9.914469363407159
$dataKL[21,0] = 23.67561901179029
$dataKL[21,1] = 9.91159716381182
$dataKL[22,0] = 22.96241447379619
$dataKL[22,1] = 9.904081729115768
$dataKL[23,0] = 22.18916664885925
$dataKL[23,1] = 9.903154333869644
$ws.Range("K2:L25").Value = $dataKL

